# Applies the "viendo solucion ejercicio 2" edit:
#  - Taylor_Orden3 (sheet 1): truncate data to rows 2-17 with a new adaptive
#    time-step sequence in column A, and zero out columns B/C.
#  - Runge_Kutta_56 (sheet 2): zero out columns B/C (column A untouched).
#  - Adams_Bashforth_Moulton (sheet 3): zero out columns B/C (column A untouched).
#  - Resumen_Comparativo (sheet 4): update function-evaluation counts and
#    error/step columns for the taylor and rk56 rows, and zero the adams
#    error column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Taylor_Orden3
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Taylor_Orden3")

$newTimes = @(
    0,
    0.0001,
    0.00022,
    0.000364,
    0.0005368,
    0.00074416,
    0.000992992,
    0.0012915904,
    0.00164990848,
    0.002079890176,
    0.0025958682112,
    0.00321504185344,
    0.003958050224128,
    0.0048496602689536,
    0.00591959232274432,
    0.006086279092414062
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $ws1.Cells.Item($i + 2, 1).Value = $newTimes[$i]
}

# Altura/Velocidad collapse to 0 for the surviving rows (2-17).
$ws1.Range("B2:C17").Value = 0

# Drop the now-unused tail (rows 18-63) so the sheet shrinks to A1:C17.
$ws1.Range("A18:A63").EntireRow.Delete()

# ---------------------------------------------------------------------------
# Sheet 2: Runge_Kutta_56
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Runge_Kutta_56")
$ws2.Range("B2:C127").Value = 0

# ---------------------------------------------------------------------------
# Sheet 3: Adams_Bashforth_Moulton
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Adams_Bashforth_Moulton")
$ws3.Range("B2:C63").Value = 0

# ---------------------------------------------------------------------------
# Sheet 4: Resumen_Comparativo
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Resumen_Comparativo")

$ws4.Range("B2").Value = 45
$ws4.Range("C2").Value = 0
$ws4.Range("D2").Value = 0.0004057519394942708

$ws4.Range("B3").Value = 77
$ws4.Range("C3").Value = 0

$ws4.Range("C4").Value = 0
